$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.474.65"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.268.93"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.66"
$ws.Range("E5").Value = "  +6.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.47"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.56"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.37"
$ws.Range("E12").Value = "  +6.78%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.82"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.920"
$ws.Range("E15").Value = "  +7.59%  "
$ws.Range("D16").Value = "2.611.03"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.270.48"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "43.605.95"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.93"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.55"
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.87"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.00"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.62"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.24"
$ws.Range("E26").Value = "  +8.52%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.76"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.80"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.54"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0920"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.33"
$ws.Range("E36").Value = "  +12.80%  "
$ws.Range("E37").Value = "  +7.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.69"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.68"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.240"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.23"
$ws.Range("E43").Value = "  -4.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.38"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.72"
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "76.35"
$ws.Range("E47").Value = "  +37.53%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.669"
$ws.Range("E48").Value = "  +19.86%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("E51").Value = "  +1.90%  "
